$d = $word.ActiveDocument

# Helper: locate a paragraph index (1-based) whose text contains the given substring.
function Get-ParaIndexByText($text) {
    $paras = $d.Paragraphs
    for ($i = 1; $i -le $paras.Count; $i++) {
        if ($paras.Item($i).Range.Text.Contains($text)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Professional summary: neutralize "all Black and Asian-American voters"
#    to "50M voters".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic demographic coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic demographic coding errors affecting 50M voters, developed",
    2
) | Out-Null

# ---------------------------------------------------------------------------
# 2) Siege Analytics bullet: same neutralization, but "50M" needs to become
#    its own bold/colored run (matching the existing 23% / 64% run style).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Discovered systematic race coding errors affecting all Black and Asian-American voters, developed",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discovered systematic race coding errors affecting 50M voters, developed",
    2
) | Out-Null

$rngBold = $d.Content
$rngBold.Find.Execute("50M", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngBold.Font.Bold = 1
$rngBold.Font.Color = 5258796   # RGB(0x2C,0x3E,0x50) packed as BGR for Word's Font.Color

# ---------------------------------------------------------------------------
# 3) Reorder PROFESSIONAL EXPERIENCE sections.
#    Before: Siege -> Mautinoa -> Salsa -> Praxis -> PCCC -> Helm/Murmuration
#    After:  Siege -> Helm/Murmuration -> Mautinoa -> PCCC -> Salsa -> Praxis
# ---------------------------------------------------------------------------

# 3a) Move the "Data Products Manager - Helm/Murmuration" block (5 paragraphs)
#     to sit right after the Siege Analytics block, before "Mautinoa".
$startIdx = Get-ParaIndexByText("Data Products Manager - Helm/Murmuration")
$endIdx = $startIdx + 4
$paras = $d.Paragraphs
$moveRng = $d.Range($paras.Item($startIdx).Range.Start, $paras.Item($endIdx).Range.End)
$moveRng.Cut() | Out-Null

$targetIdx = Get-ParaIndexByText("Software Engineer - Mautinoa Technologies")
$paras = $d.Paragraphs
$target = $paras.Item($targetIdx)
$pasteRng = $d.Range($target.Range.Start, $target.Range.Start)
$pasteRng.Paste() | Out-Null

# Restore the Heading3 style on the pasted "Data Products Manager" heading
# (Paste collapses into the destination paragraph's formatting).
$paras = $d.Paragraphs
$helmHeadingIdx = Get-ParaIndexByText("Data Products Manager - Helm/Murmuration")
$paras.Item($helmHeadingIdx).Style = "Heading 3"

# Restore bold/colored "57%" run lost during the paste.
$rngPct = $d.Content
$rngPct.Find.Execute("57%", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngPct.Font.Bold = 1
$rngPct.Font.Color = 5258796

# 3b) Move the "Research Director - PCCC" block (5 paragraphs) to sit right
#     after the Mautinoa block, before "Salsa Labs".
$startIdx2 = Get-ParaIndexByText("Research Director - PCCC")
$endIdx2 = $startIdx2 + 4
$paras = $d.Paragraphs
$moveRng2 = $d.Range($paras.Item($startIdx2).Range.Start, $paras.Item($endIdx2).Range.End)
$moveRng2.Cut() | Out-Null

$targetIdx2 = Get-ParaIndexByText("Software Engineer - Salsa Labs")
$paras = $d.Paragraphs
$target2 = $paras.Item($targetIdx2)
$pasteRng2 = $d.Range($target2.Range.Start, $target2.Range.Start)
$pasteRng2.Paste() | Out-Null

# Restore the Heading3 style on the pasted "Research Director" heading.
$paras = $d.Paragraphs
$pcccHeadingIdx = Get-ParaIndexByText("Research Director - PCCC")
$paras.Item($pcccHeadingIdx).Style = "Heading 3"

Write-Output "done"
